$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text renames ---
# (order matters for the resulting shared-string table layout)
$ws.Range("E2").Value = "Átlagos feldolgozási idő (s)"
$ws.Range("D2").Value = "Pontosság (%)"
$ws.Range("B2").Value = "LLM modell"

# --- Header row (B2:H2): vertical alignment top -> center ---
$ws.Range("B2:H2").VerticalAlignment = -4108

# --- Data formatting cleanup ---
# Bring column B (and D:H, which previously had no alignment set) in line
# with column C's existing "centered, bordered" look, by copying the
# already-correct format from column C onto the rest of each row. This
# reuses the existing style entries instead of synthesizing new ones.

# Rows 3, 9, 15 ("first of group" rows -> yellow-filled style, like C3)
$ws.Range("C3").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null
$ws.Range("D3:H3").PasteSpecial(-4122) | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null
$ws.Range("D9:H9").PasteSpecial(-4122) | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$ws.Range("D15:H15").PasteSpecial(-4122) | Out-Null

# Remaining rows (no-fill style, like C4)
$ws.Range("C4").Copy() | Out-Null
$ws.Range("B4:B8").PasteSpecial(-4122) | Out-Null
$ws.Range("D4:H8").PasteSpecial(-4122) | Out-Null
$ws.Range("B10:B14").PasteSpecial(-4122) | Out-Null
$ws.Range("D10:H14").PasteSpecial(-4122) | Out-Null
$ws.Range("B16:B20").PasteSpecial(-4122) | Out-Null
$ws.Range("D16:H20").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Column width changes ---
$ws.Columns("F").ColumnWidth = 18
$ws.Columns("G").ColumnWidth = 12.666666666666666

# --- Cursor / selection ---
$ws.Range("C24").Select() | Out-Null
